# Adds example motor/propeller performance rows ("Made some examples for UMD")
# to the "List" summary sheet and detailed per-throttle data to the "KD" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "KD" sheet — detailed throttle/power/thrust/RPM tables for the new
#    KDE2315XF-885 motor across 11.6", 15.4" and 23.1" voltage classes and
#    five new propeller sizes.
# ---------------------------------------------------------------------------
$wsKD = $wb.Worksheets.Item("KD")

$kdRows = @(
    @(1031, "11.6", "9” x 3.0", 0.25, 6, 80, 2820),
    @(1032, $null, $null, 0.375, 13, 150, 4280),
    @(1033, $null, $null, 0.5, 22, 230, 5520),
    @(1034, $null, $null, 0.625, 35, 330, 6660),
    @(1035, $null, $null, 0.75, 54, 450, 7860),
    @(1036, $null, $null, 0.875, 85, 610, 9060),
    @(1037, $null, $null, 1, 113, 740, 9960),
    @(1038, "11.6", "9`" x 4.5 (DJI)", 0.25, 7, 110, 2700),
    @(1039, $null, $null, 0.375, 15, 190, 3900),
    @(1040, $null, $null, 0.5, 33, 320, 5040),
    @(1041, $null, $null, 0.625, 58, 480, 6180),
    @(1042, $null, $null, 0.75, 90, 640, 7200),
    @(1043, $null, $null, 0.875, 135, 820, 8160),
    @(1044, $null, $null, 1, 174, 970, 8940),
    @(1045, "11.6", "10`" x 3.3", 0.25, 7, 110, 2720),
    @(1046, $null, $null, 0.375, 15, 200, 4140),
    @(1047, $null, $null, 0.5, 29, 300, 5340),
    @(1048, $null, $null, 0.625, 48, 430, 6480),
    @(1049, $null, $null, 0.75, 74, 590, 7500),
    @(1050, $null, $null, 0.875, 112, 770, 8580),
    @(1051, $null, $null, 1, 148, 940, 9480),
    @(1052, "11.6", "11`" x 3.7", 0.25, 7, 120, 2640),
    @(1053, $null, $null, 0.375, 19, 240, 3840),
    @(1054, $null, $null, 0.5, 38, 400, 4980),
    @(1055, $null, $null, 0.625, 66, 570, 6000),
    @(1056, $null, $null, 0.75, 102, 760, 6960),
    @(1057, $null, $null, 0.875, 150, 990, 7860),
    @(1058, $null, $null, 1, 196, 1170, 8580),
    @(1059, "11.6", "12” x 4.0", 0.25, 10, 160, 2520),
    @(1060, $null, $null, 0.375, 22, 280, 3600),
    @(1061, $null, $null, 0.5, 48, 480, 4560),
    @(1062, $null, $null, 0.625, 83, 670, 5520),
    @(1063, $null, $null, 0.75, 130, 880, 6300),
    @(1064, $null, $null, 0.875, 190, 1100, 7020),
    @(1065, $null, $null, 1, 243, 1280, 7560),
    @(1066, "15.4", "9” x 3.0", 0.25, 10, 130, 3720),
    @(1067, $null, $null, 0.375, 23, 230, 5460),
    @(1068, $null, $null, 0.5, 44, 370, 7140),
    @(1069, $null, $null, 0.625, 72, 540, 8520),
    @(1070, $null, $null, 0.75, 113, 740, 9980),
    @(1071, $null, $null, 0.875, 178, 980, 11480),
    @(1072, $null, $null, 1, 232, 1160, 12420),
    @(1073, "15.4", "9`" x 4.5 (DJI)", 0.25, 13, 180, 3480),
    @(1074, $null, $null, 0.375, 33, 320, 5040),
    @(1075, $null, $null, 0.5, 67, 510, 6420),
    @(1076, $null, $null, 0.625, 113, 730, 7740),
    @(1077, $null, $null, 0.75, 176, 970, 8910),
    @(1078, $null, $null, 0.875, 259, 1230, 9960),
    @(1079, $null, $null, 1, 335, 1440, 10920),
    @(1080, "15.4", "10`" x 3.3", 0.25, 11, 170, 3660),
    @(1081, $null, $null, 0.375, 28, 300, 5220),
    @(1082, $null, $null, 0.5, 55, 470, 6720),
    @(1083, $null, $null, 0.625, 95, 690, 8180),
    @(1084, $null, $null, 0.75, 147, 920, 9360),
    @(1085, $null, $null, 0.875, 222, 1210, 10680),
    @(1086, $null, $null, 1, 290, 1430, 11640),
    @(1087, "15.4", "11`" x 3.7", 0.25, 13, 210, 3420),
    @(1088, $null, $null, 0.375, 37, 390, 4860),
    @(1089, $null, $null, 0.5, 74, 610, 6180),
    @(1090, $null, $null, 0.625, 127, 870, 7380),
    @(1091, $null, $null, 0.75, 198, 1160, 8520),
    @(1092, $null, $null, 0.875, 290, 1440, 9480),
    @(1093, $null, $null, 1, 368, 1650, 10320),
    @(1094, "15.4", "12” x 4.0", 0.25, 16, 250, 3180),
    @(1095, $null, $null, 0.375, 45, 440, 4440),
    @(1096, $null, $null, 0.5, 96, 720, 5700),
    @(1097, $null, $null, 0.625, 161, 980, 6660),
    @(1098, $null, $null, 0.75, 242, 1260, 7500),
    @(1099, $null, $null, 0.875, 344, 1500, 8160),
    @(1100, $null, $null, 1, 456, 1820, 9040),
    @(1101, "23.1", "9” x 3.0", 0.25, 25, 230, 5460),
    @(1102, $null, $null, 0.375, 63, 450, 7800),
    @(1103, $null, $null, 0.5, 119, 720, 9720),
    @(1104, $null, $null, 0.625, 198, 1030, 11640),
    @(1105, $null, $null, 0.75, 318, 1390, 13440),
    @(1106, $null, $null, 0.875, 478, 1760, 15060),
    @(1107, $null, $null, 1, 605, 2020, 16260)
)

foreach ($row in $kdRows) {
    $r = $row[0]
    if ($row[1] -ne $null) {
        $wsKD.Cells.Item($r, 1).Value = "KDE2315XF-885"
        $wsKD.Cells.Item($r, 2).Value = [double]$row[1]
        $wsKD.Cells.Item($r, 3).Value = $row[2]
    }
    $wsKD.Cells.Item($r, 4).Value = [double]$row[3]
    $wsKD.Cells.Item($r, 4).NumberFormat = "0.00%"
    $wsKD.Cells.Item($r, 5).Value = [double]$row[4]
    $wsKD.Cells.Item($r, 6).Value = [double]$row[5]
    $wsKD.Cells.Item($r, 7).Value = [double]$row[6]
}

# ---------------------------------------------------------------------------
# 2. "List" sheet — summary rows for the new KDE2315XF-885 motor paired with
#    its five new propeller options. Mass (column E) is computed with the
#    same "=75+10" formula used by the sibling example rows.
# ---------------------------------------------------------------------------
$wsList = $wb.Worksheets.Item("List")

$listRows = @(
    @(118, "KDE2315XF-885", "9” x 3.0", 9, 17000, 24),
    @(119, "KDE2315XF-885", "9`" x 4.5 (DJI)", 9, 17000, 24),
    @(120, "KDE2315XF-885", "10`" x 3.3", 10, 12000, 24),
    @(121, "KDE2315XF-885", "11`" x 3.7", 11, 11000, 24),
    @(122, "KDE2315XF-885", "12” x 4.0", 12, 10000, 24)
)

foreach ($row in $listRows) {
    $r = $row[0]
    $wsList.Cells.Item($r, 1).Value = $row[1]
    $wsList.Cells.Item($r, 2).Value = $row[2]
    $wsList.Cells.Item($r, 3).Value = [double]$row[3]
    $wsList.Cells.Item($r, 4).Value = [double]$row[4]
    $wsList.Cells.Item($r, 5).Formula = "=75+10"
    $wsList.Cells.Item($r, 6).Value = [double]$row[5]
    $wsList.Cells.Item($r, 7).Value = 1
}

# ---------------------------------------------------------------------------
# 3. Leave the UI focused on the last row typed, on the "List" sheet.
# ---------------------------------------------------------------------------
$wsList.Activate()
$wsList.Range("G122").Select()
